$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: bump "Forandrad" (column C) date for all data rows 2-49 to 2026-02-19 (serial 46072) ---
$ws.Range("C2:C49").Value = 46072

# --- Step 2: re-write the rows whose entire content was permuted/reordered (rows 9,10 / 14,15 / 28-38) ---

# Row 9 (content formerly at row 10)
$ws.Range("A9").Value = 'A 4792-2025'
$ws.Range("B9").Value = 45688.57549768518
$ws.Range("D9").Value = 'SKÅNE LÄN'
$ws.Range("E9").Value = 'LUND'
$ws.Range("F9").ClearContents()
$ws.Range("G9").Value = 5.9
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 1
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = 0
$ws.Range("O9").Value = 1
$ws.Range("P9").Value = 0
$ws.Range("Q9").Value = 1
$ws.Range("R9").Value = 'Oxtungssvamp'
$ws.Range("S9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/artfynd/A 4792-2025 artfynd.xlsx", "A 4792-2025")'
$ws.Range("T9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/kartor/A 4792-2025 karta.png", "A 4792-2025")'
$ws.Range("V9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/klagomål/A 4792-2025 FSC-klagomål.docx", "A 4792-2025")'
$ws.Range("W9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/klagomålsmail/A 4792-2025 FSC-klagomål mail.docx", "A 4792-2025")'
$ws.Range("X9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/tillsyn/A 4792-2025 tillsynsbegäran.docx", "A 4792-2025")'
$ws.Range("Y9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/tillsynsmail/A 4792-2025 tillsynsbegäran mail.docx", "A 4792-2025")'

# Row 10 (content formerly at row 9)
$ws.Range("A10").Value = 'A 11732-2024'
$ws.Range("B10").Value = 45373.69222222222
$ws.Range("D10").Value = 'SKÅNE LÄN'
$ws.Range("E10").Value = 'LUND'
$ws.Range("F10").ClearContents()
$ws.Range("G10").Value = 2.8
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 0
$ws.Range("O10").Value = 1
$ws.Range("P10").Value = 1
$ws.Range("Q10").Value = 1
$ws.Range("R10").Value = 'Dvärgjohannesört'
$ws.Range("S10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/artfynd/A 11732-2024 artfynd.xlsx", "A 11732-2024")'
$ws.Range("T10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/kartor/A 11732-2024 karta.png", "A 11732-2024")'
$ws.Range("V10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/klagomål/A 11732-2024 FSC-klagomål.docx", "A 11732-2024")'
$ws.Range("W10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/klagomålsmail/A 11732-2024 FSC-klagomål mail.docx", "A 11732-2024")'
$ws.Range("X10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/tillsyn/A 11732-2024 tillsynsbegäran.docx", "A 11732-2024")'
$ws.Range("Y10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/tillsynsmail/A 11732-2024 tillsynsbegäran mail.docx", "A 11732-2024")'

# Row 14 (content formerly at row 15)
$ws.Range("A14").Value = 'A 59227-2025'
$ws.Range("B14").Value = 45988.62253472222
$ws.Range("D14").Value = 'SKÅNE LÄN'
$ws.Range("E14").Value = 'LUND'
$ws.Range("F14").ClearContents()
$ws.Range("G14").Value = 7.1
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 1
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = 0
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = 0
$ws.Range("Q14").Value = 1
$ws.Range("R14").Value = 'Scharlakansskål'
$ws.Range("S14").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/artfynd/A 59227-2025 artfynd.xlsx", "A 59227-2025")'
$ws.Range("T14").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/kartor/A 59227-2025 karta.png", "A 59227-2025")'
$ws.Range("V14").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/klagomål/A 59227-2025 FSC-klagomål.docx", "A 59227-2025")'
$ws.Range("W14").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/klagomålsmail/A 59227-2025 FSC-klagomål mail.docx", "A 59227-2025")'
$ws.Range("X14").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/tillsyn/A 59227-2025 tillsynsbegäran.docx", "A 59227-2025")'
$ws.Range("Y14").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/tillsynsmail/A 59227-2025 tillsynsbegäran mail.docx", "A 59227-2025")'

# Row 15 (content formerly at row 14)
$ws.Range("A15").Value = 'A 35443-2021'
$ws.Range("B15").Value = 44385
$ws.Range("D15").Value = 'SKÅNE LÄN'
$ws.Range("E15").Value = 'LUND'
$ws.Range("F15").ClearContents()
$ws.Range("G15").Value = 4.2
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = 0
$ws.Range("O15").Value = 1
$ws.Range("P15").Value = 1
$ws.Range("Q15").Value = 1
$ws.Range("R15").Value = 'Klubbfibbla'
$ws.Range("S15").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/artfynd/A 35443-2021 artfynd.xlsx", "A 35443-2021")'
$ws.Range("T15").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/kartor/A 35443-2021 karta.png", "A 35443-2021")'
$ws.Range("V15").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/klagomål/A 35443-2021 FSC-klagomål.docx", "A 35443-2021")'
$ws.Range("W15").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/klagomålsmail/A 35443-2021 FSC-klagomål mail.docx", "A 35443-2021")'
$ws.Range("X15").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/tillsyn/A 35443-2021 tillsynsbegäran.docx", "A 35443-2021")'
$ws.Range("Y15").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/tillsynsmail/A 35443-2021 tillsynsbegäran mail.docx", "A 35443-2021")'

# Row 28 (content formerly at row 34)
$ws.Range("A28").Value = 'A 58335-2025'
$ws.Range("B28").Value = 45985.48332175926
$ws.Range("D28").Value = 'SKÅNE LÄN'
$ws.Range("E28").Value = 'LUND'
$ws.Range("F28").ClearContents()
$ws.Range("G28").Value = 2.5
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 0
$ws.Range("N28").Value = 0
$ws.Range("O28").Value = 0
$ws.Range("P28").Value = 0
$ws.Range("Q28").Value = 0
$ws.Range("R28").ClearContents()
$ws.Range("S28").ClearContents()
$ws.Range("T28").ClearContents()
$ws.Range("V28").ClearContents()
$ws.Range("W28").ClearContents()
$ws.Range("X28").ClearContents()
$ws.Range("Y28").ClearContents()

# Row 29 (content formerly at row 35)
$ws.Range("A29").Value = 'A 58488-2025'
$ws.Range("B29").Value = 45985.66210648148
$ws.Range("D29").Value = 'SKÅNE LÄN'
$ws.Range("E29").Value = 'LUND'
$ws.Range("F29").ClearContents()
$ws.Range("G29").Value = 9.5
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 0
$ws.Range("N29").Value = 0
$ws.Range("O29").Value = 0
$ws.Range("P29").Value = 0
$ws.Range("Q29").Value = 0
$ws.Range("R29").ClearContents()
$ws.Range("S29").ClearContents()
$ws.Range("T29").ClearContents()
$ws.Range("V29").ClearContents()
$ws.Range("W29").ClearContents()
$ws.Range("X29").ClearContents()
$ws.Range("Y29").ClearContents()

# Row 30 (content formerly at row 37)
$ws.Range("A30").Value = 'A 58337-2025'
$ws.Range("B30").Value = 45985.48895833334
$ws.Range("D30").Value = 'SKÅNE LÄN'
$ws.Range("E30").Value = 'LUND'
$ws.Range("F30").ClearContents()
$ws.Range("G30").Value = 2
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = 0
$ws.Range("N30").Value = 0
$ws.Range("O30").Value = 0
$ws.Range("P30").Value = 0
$ws.Range("Q30").Value = 0
$ws.Range("R30").ClearContents()
$ws.Range("S30").ClearContents()
$ws.Range("T30").ClearContents()
$ws.Range("V30").ClearContents()
$ws.Range("W30").ClearContents()
$ws.Range("X30").ClearContents()
$ws.Range("Y30").ClearContents()

# Row 31 (content formerly at row 28)
$ws.Range("A31").Value = 'A 29247-2025'
$ws.Range("B31").Value = 45824
$ws.Range("D31").Value = 'SKÅNE LÄN'
$ws.Range("E31").Value = 'LUND'
$ws.Range("F31").ClearContents()
$ws.Range("G31").Value = 0.9
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = 0
$ws.Range("N31").Value = 0
$ws.Range("O31").Value = 0
$ws.Range("P31").Value = 0
$ws.Range("Q31").Value = 0
$ws.Range("R31").ClearContents()
$ws.Range("S31").ClearContents()
$ws.Range("T31").ClearContents()
$ws.Range("V31").ClearContents()
$ws.Range("W31").ClearContents()
$ws.Range("X31").ClearContents()
$ws.Range("Y31").ClearContents()

# Row 32 (content formerly at row 30)
$ws.Range("A32").Value = 'A 28885-2022'
$ws.Range("B32").Value = 44749
$ws.Range("D32").Value = 'SKÅNE LÄN'
$ws.Range("E32").Value = 'LUND'
$ws.Range("F32").ClearContents()
$ws.Range("G32").Value = 2.4
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = 0
$ws.Range("N32").Value = 0
$ws.Range("O32").Value = 0
$ws.Range("P32").Value = 0
$ws.Range("Q32").Value = 0
$ws.Range("R32").ClearContents()
$ws.Range("S32").ClearContents()
$ws.Range("T32").ClearContents()
$ws.Range("V32").ClearContents()
$ws.Range("W32").ClearContents()
$ws.Range("X32").ClearContents()
$ws.Range("Y32").ClearContents()

# Row 33 (content formerly at row 29)
$ws.Range("A33").Value = 'A 34484-2025'
$ws.Range("B33").Value = 45847.39517361111
$ws.Range("D33").Value = 'SKÅNE LÄN'
$ws.Range("E33").Value = 'LUND'
$ws.Range("F33").Value = 'Kommuner'
$ws.Range("G33").Value = 1.2
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = 0
$ws.Range("N33").Value = 0
$ws.Range("O33").Value = 0
$ws.Range("P33").Value = 0
$ws.Range("Q33").Value = 0
$ws.Range("R33").ClearContents()
$ws.Range("S33").ClearContents()
$ws.Range("T33").ClearContents()
$ws.Range("V33").ClearContents()
$ws.Range("W33").ClearContents()
$ws.Range("X33").ClearContents()
$ws.Range("Y33").ClearContents()

# Row 34 (content formerly at row 31)
$ws.Range("A34").Value = 'A 32972-2024'
$ws.Range("B34").Value = 45517.42064814815
$ws.Range("D34").Value = 'SKÅNE LÄN'
$ws.Range("E34").Value = 'LUND'
$ws.Range("F34").Value = 'Kommuner'
$ws.Range("G34").Value = 20.7
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = 0
$ws.Range("N34").Value = 0
$ws.Range("O34").Value = 0
$ws.Range("P34").Value = 0
$ws.Range("Q34").Value = 0
$ws.Range("R34").ClearContents()
$ws.Range("S34").ClearContents()
$ws.Range("T34").ClearContents()
$ws.Range("V34").ClearContents()
$ws.Range("W34").ClearContents()
$ws.Range("X34").ClearContents()
$ws.Range("Y34").ClearContents()

# Row 35 (content formerly at row 38)
$ws.Range("A35").Value = 'A 59206-2025'
$ws.Range("B35").Value = 45988.60881944445
$ws.Range("D35").Value = 'SKÅNE LÄN'
$ws.Range("E35").Value = 'LUND'
$ws.Range("F35").ClearContents()
$ws.Range("G35").Value = 6.9
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = 0
$ws.Range("N35").Value = 0
$ws.Range("O35").Value = 0
$ws.Range("P35").Value = 0
$ws.Range("Q35").Value = 0
$ws.Range("R35").ClearContents()
$ws.Range("S35").ClearContents()
$ws.Range("T35").ClearContents()
$ws.Range("V35").ClearContents()
$ws.Range("W35").ClearContents()
$ws.Range("X35").ClearContents()
$ws.Range("Y35").ClearContents()

# Row 36 (content formerly at row 32)
$ws.Range("A36").Value = 'A 29075-2023'
$ws.Range("B36").Value = 45104
$ws.Range("D36").Value = 'SKÅNE LÄN'
$ws.Range("E36").Value = 'LUND'
$ws.Range("F36").ClearContents()
$ws.Range("G36").Value = 0.7
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = 0
$ws.Range("N36").Value = 0
$ws.Range("O36").Value = 0
$ws.Range("P36").Value = 0
$ws.Range("Q36").Value = 0
$ws.Range("R36").ClearContents()
$ws.Range("S36").ClearContents()
$ws.Range("T36").ClearContents()
$ws.Range("V36").ClearContents()
$ws.Range("W36").ClearContents()
$ws.Range("X36").ClearContents()
$ws.Range("Y36").ClearContents()

# Row 37 (content formerly at row 33)
$ws.Range("A37").Value = 'A 33985-2021'
$ws.Range("B37").Value = 44378
$ws.Range("D37").Value = 'SKÅNE LÄN'
$ws.Range("E37").Value = 'LUND'
$ws.Range("F37").ClearContents()
$ws.Range("G37").Value = 1.5
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = 0
$ws.Range("N37").Value = 0
$ws.Range("O37").Value = 0
$ws.Range("P37").Value = 0
$ws.Range("Q37").Value = 0
$ws.Range("R37").ClearContents()
$ws.Range("S37").ClearContents()
$ws.Range("T37").ClearContents()
$ws.Range("V37").ClearContents()
$ws.Range("W37").ClearContents()
$ws.Range("X37").ClearContents()
$ws.Range("Y37").ClearContents()

# Row 38 (content formerly at row 36)
$ws.Range("A38").Value = 'A 45725-2022'
$ws.Range("B38").Value = 44846
$ws.Range("D38").Value = 'SKÅNE LÄN'
$ws.Range("E38").Value = 'LUND'
$ws.Range("F38").ClearContents()
$ws.Range("G38").Value = 1.2
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 0
$ws.Range("N38").Value = 0
$ws.Range("O38").Value = 0
$ws.Range("P38").Value = 0
$ws.Range("Q38").Value = 0
$ws.Range("R38").ClearContents()
$ws.Range("S38").ClearContents()
$ws.Range("T38").ClearContents()
$ws.Range("V38").ClearContents()
$ws.Range("W38").ClearContents()
$ws.Range("X38").ClearContents()
$ws.Range("Y38").ClearContents()
